{"js": "// Update the date/weekday line at the top of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-10-19 Saturday\", Word.InsertLocation.replace);\n\n// Update the five populated rows of the multiplication-practice table.\n// rowIndex is 0-based against the full table grid (20 rows x 5 cols);\n// only rows 0, 4, 9, 14 and 19 contain text in this document.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst updates = [\n  [0, [\"60\u00d767=\", \"86\u00d791=\", \"74\u00d740=\", \"96\u00d781=\", \"92\u00d711=\"]],\n  [4, [\"61\u00d791=\", \"94\u00d740=\", \"50\u00d711=\", \"73\u00d761=\", \"49\u00d755=\"]],\n  [9, [\"44\u00d715=\", \"31\u00d746=\", \"74\u00d778=\", \"31\u00d771=\", \"74\u00d735=\"]],\n  [14, [\"49\u00d792=\", \"44\u00d791=\", \"88\u00d729=\", \"87\u00d735=\", \"64\u00d790=\"]],\n  [19, [\"34\u00d771=\", \"44\u00d794=\", \"32\u00d793=\", \"98\u00d730=\", \"84\u00d793=\"]],\n];\n\nfor (const [rowIndex, values] of updates) {\n  for (let colIndex = 0; colIndex < values.length; colIndex++) {\n    const cell = table.getCellOrNullObject(rowIndex, colIndex);\n    cell.value = values[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date/weekday line at the top of the document.\n$d.Paragraphs.Item(1).Range.Text = \"2024-10-19 Saturday\"\n\n# Update the five populated rows of the multiplication-practice table.\n# Each inner array is the new text for columns 1..5 of that table row.\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = @(\"60\u00d767=\", \"86\u00d791=\", \"74\u00d740=\", \"96\u00d781=\", \"92\u00d711=\")\n    5  = @(\"61\u00d791=\", \"94\u00d740=\", \"50\u00d711=\", \"73\u00d761=\", \"49\u00d755=\")\n    10 = @(\"44\u00d715=\", \"31\u00d746=\", \"74\u00d778=\", \"31\u00d771=\", \"74\u00d735=\")\n    15 = @(\"49\u00d792=\", \"44\u00d791=\", \"88\u00d729=\", \"87\u00d735=\", \"64\u00d790=\")\n    20 = @(\"34\u00d771=\", \"44\u00d794=\", \"32\u00d793=\", \"98\u00d730=\", \"84\u00d793=\")\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $values = $updates[$rowIndex]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
